$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.9448398576512456
$ws1.Range("C2").Value = 0.4545454545454545
$ws1.Range("D2").Value = 0.5357142857142857
$ws1.Range("E2").Value = 0.4918032786885246
$ws1.Range("F2").Value = 0.5172413793103449
$ws1.Range("G2").Value = 0.5320600272851296
$ws1.Range("H2").Value = 0.7510032102728731
$ws1.Range("I2").Value = 15
$ws1.Range("J2").Value = 18
$ws1.Range("K2").Value = 516
$ws1.Range("L2").Value = 13

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 - class "0"
$ws2.Range("B2").Value = 0.9754253308128544
$ws2.Range("C2").Value = 0.9662921348314607
$ws2.Range("D2").Value = 0.9708372530573848

# row 3 - class "1"
$ws2.Range("B3").Value = 0.4545454545454545
$ws2.Range("C3").Value = 0.5357142857142857
$ws2.Range("D3").Value = 0.4918032786885246

# row 4 - accuracy
$ws2.Range("B4").Value = 0.9448398576512456
$ws2.Range("C4").Value = 0.9448398576512456
$ws2.Range("D4").Value = 0.9448398576512456
$ws2.Range("E4").Value = 0.9448398576512456

# row 5 - macro avg
$ws2.Range("B5").Value = 0.7149853926791545
$ws2.Range("C5").Value = 0.7510032102728732
$ws2.Range("D5").Value = 0.7313202658729547

# row 6 - weighted avg
$ws2.Range("B6").Value = 0.9494740202514893
$ws2.Range("C6").Value = 0.9448398576512456
$ws2.Range("D6").Value = 0.946970791700929

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 516
$ws3.Range("C2").Value = 18
$ws3.Range("B3").Value = 13
$ws3.Range("C3").Value = 15
